# Applies the "Lab 3 Inverse Kinematics" report revision:
#  - clarifies the (0,0) target point to (0,0,0)
#  - adds radius/units context to the results table header
#  - fills in the previously-empty "Measured Error" cells
#  - expands the reflection paragraphs with the atan2/heading fix detail
#  - notes that GitHub Copilot was used through VS Code

$d = $word.ActiveDocument

# 1. "(0,0)" -> "(0,0,0)" in the intro paragraph about measured error distance.
$d.Content.Find.Execute("(0,0)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(0,0,0)", 2) | Out-Null

$tbl = $d.Tables(1)

# 2. Header row, first cell: append ", radius = 50 cm" after "Waypoint Navigation Run".
$headCell1 = $tbl.Cell(1, 1)
$headCell1.Range.Paragraphs(1).Range.InsertAfter(", radius = 50 cm")

# 3. Header row, second cell: "Measured Error (cm)" -> "Measured Error (cm, cm, deg)".
$d.Content.Find.Execute("Measured Error (cm)", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Measured Error (cm, cm, deg)", 2) | Out-Null

# 4. Fill in the three previously empty "Measured Error" data cells.
$tbl.Cell(2, 2).Range.Text = "(-9, 0, 0 deg)"
$tbl.Cell(3, 2).Range.Text = "(9, 3, 13 deg )"
$tbl.Cell(4, 2).Range.Text = "(6, 2, 4.5 deg)"

# 5. Expand the "most time spent" reflection paragraph.
$oldReflection = "This was definitely the lab I spent the most time on. There was " + `
    "significant debugging to be done with the waypoint algorithm. I spent maybe " + `
    "8-10 hours on it all told, a certain portion of that on getting proper data " + `
    "reporting on the robot. "
$newReflection = "This was definitely the lab I spent the most time on so far. There was " + `
    "significant debugging to be done with the waypoint algorithm. I spent maybe 15 " + `
    "hours on it all told, a certain portion of that on getting proper data reporting " + `
    "on the robot, and a major portion of that struggling to realize that I needed to " + `
    "subtract initial heading from the atan2 output to get the proper alpha turn. "
$d.Content.Find.Execute($oldReflection, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newReflection, 2) | Out-Null

# 6. Note that Copilot was used through VS Code.
$d.Content.Find.Execute("I used GitHub Copilot on some", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "I used GitHub Copilot through VS Code on some", 2) | Out-Null
